# Protocol plan.xlsx - "Added in synchronisation for units."
# Adds a new "country, money" column (J) to the protocol table on Sheet1,
# tweaks a couple of existing cell values, and moves the table's
# selection/right-hand border accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell value edits ---------------------------------------------------

# Parametres row: "id, x, y" -> " x, y, id"
$ws.Range("D4").Value = " x, y, id"

# Returns row, column I: "Country number" -> "country"
$ws.Range("I4").Value = "country"

# New column J content
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = "Money"
$ws.Range("J4").Value = "country, money"
$ws.Range("J5").Value = "-"
$ws.Range("J6").Value = "-"

# --- Column width ---------------------------------------------------------
# New column J, matches the other data columns (stored width 15)
$ws.Columns.Item(10).ColumnWidth = 14.14

# --- Borders ----------------------------------------------------------
# Column I used to be the right-hand edge of the table (medium right
# border). Column J now takes over as the right edge: it gets the medium
# right border plus whatever top/bottom edge column I had, while column
# I's own right edge becomes borderless (interior column gap).

$xlContinuous = 1
$xlThin = 2
$xlMedium = -4138
$xlNone = -4142

# Column I: drop the old medium right border on every row.
for ($row = 1; $row -le 6; $row++) {
    $ws.Cells.Item($row, 9).Borders.Item(10).LineStyle = $xlNone
}
# Row 1 additionally gains a thin bottom border on column I (matching the
# rest of the header row, e.g. E1:H1).
$ws.Cells.Item(1, 9).Borders.Item(9).LineStyle = $xlContinuous
$ws.Cells.Item(1, 9).Borders.Item(9).Weight = $xlThin

# Column J: new right-hand table edge (medium), plus top/bottom to match
# the row it's in.
$ws.Cells.Item(1, 10).Borders.Item(10).LineStyle = $xlContinuous
$ws.Cells.Item(1, 10).Borders.Item(10).Weight = $xlMedium
$ws.Cells.Item(1, 10).Borders.Item(8).LineStyle = $xlContinuous
$ws.Cells.Item(1, 10).Borders.Item(8).Weight = $xlMedium
$ws.Cells.Item(1, 10).Borders.Item(9).LineStyle = $xlContinuous
$ws.Cells.Item(1, 10).Borders.Item(9).Weight = $xlThin

$ws.Cells.Item(2, 10).Borders.Item(10).LineStyle = $xlContinuous
$ws.Cells.Item(2, 10).Borders.Item(10).Weight = $xlMedium
$ws.Cells.Item(2, 10).Borders.Item(8).LineStyle = $xlContinuous
$ws.Cells.Item(2, 10).Borders.Item(8).Weight = $xlThin
$ws.Cells.Item(2, 10).Borders.Item(9).LineStyle = $xlContinuous
$ws.Cells.Item(2, 10).Borders.Item(9).Weight = $xlThin

$ws.Cells.Item(3, 10).Borders.Item(10).LineStyle = $xlContinuous
$ws.Cells.Item(3, 10).Borders.Item(10).Weight = $xlMedium

$ws.Cells.Item(4, 10).Borders.Item(10).LineStyle = $xlContinuous
$ws.Cells.Item(4, 10).Borders.Item(10).Weight = $xlMedium

$ws.Cells.Item(5, 10).Borders.Item(10).LineStyle = $xlContinuous
$ws.Cells.Item(5, 10).Borders.Item(10).Weight = $xlMedium

$ws.Cells.Item(6, 10).Borders.Item(10).LineStyle = $xlContinuous
$ws.Cells.Item(6, 10).Borders.Item(10).Weight = $xlMedium
$ws.Cells.Item(6, 10).Borders.Item(9).LineStyle = $xlContinuous
$ws.Cells.Item(6, 10).Borders.Item(9).Weight = $xlMedium

# --- Selection ------------------------------------------------------------
$ws.Range("F5").Select()
